$d = $word.ActiveDocument
$d.Content.Find.Execute("69×86=5934", $true, $false, $false, $false, $false, $true, 1, $false, "67×84=5628", 2) | Out-Null
$d.Content.Find.Execute("23×50=1150", $true, $false, $false, $false, $false, $true, 1, $false, "43×73=3139", 2) | Out-Null
$d.Content.Find.Execute("23×33=759", $true, $false, $false, $false, $false, $true, 1, $false, "71×21=1491", 2) | Out-Null
$d.Content.Find.Execute("91×79=7189", $true, $false, $false, $false, $false, $true, 1, $false, "55×75=4125", 2) | Out-Null
$d.Content.Find.Execute("82×22=1804", $true, $false, $false, $false, $false, $true, 1, $false, "64×74=4736", 2) | Out-Null
$d.Content.Find.Execute("96×24=2304", $true, $false, $false, $false, $false, $true, 1, $false, "81×76=6156", 2) | Out-Null
$d.Content.Find.Execute("31×44=1364", $true, $false, $false, $false, $false, $true, 1, $false, "53×97=5141", 2) | Out-Null
$d.Content.Find.Execute("72×42=3024", $true, $false, $false, $false, $false, $true, 1, $false, "21×13=273", 2) | Out-Null
$d.Content.Find.Execute("17×37=629", $true, $false, $false, $false, $false, $true, 1, $false, "54×56=3024", 2) | Out-Null
$d.Content.Find.Execute("54×98=5292", $true, $false, $false, $false, $false, $true, 1, $false, "35×36=1260", 2) | Out-Null
$d.Content.Find.Execute("65×17=1105", $true, $false, $false, $false, $false, $true, 1, $false, "23×15=345", 2) | Out-Null
$d.Content.Find.Execute("18×96=1728", $true, $false, $false, $false, $false, $true, 1, $false, "61×22=1342", 2) | Out-Null
$d.Content.Find.Execute("63×36=2268", $true, $false, $false, $false, $false, $true, 1, $false, "74×11=814", 2) | Out-Null
$d.Content.Find.Execute("88×12=1056", $true, $false, $false, $false, $false, $true, 1, $false, "49×28=1372", 2) | Out-Null
$d.Content.Find.Execute("14×35=490", $true, $false, $false, $false, $false, $true, 1, $false, "66×55=3630", 2) | Out-Null
$d.Content.Find.Execute("95×67=6365", $true, $false, $false, $false, $false, $true, 1, $false, "56×14=784", 2) | Out-Null
$d.Content.Find.Execute("37×18=666", $true, $false, $false, $false, $false, $true, 1, $false, "42×95=3990", 2) | Out-Null
$d.Content.Find.Execute("89×42=3738", $true, $false, $false, $false, $false, $true, 1, $false, "47×82=3854", 2) | Out-Null
$d.Content.Find.Execute("31×84=2604", $true, $false, $false, $false, $false, $true, 1, $false, "90×23=2070", 2) | Out-Null
$d.Content.Find.Execute("34×98=3332", $true, $false, $false, $false, $false, $true, 1, $false, "37×29=1073", 2) | Out-Null
$d.Content.Find.Execute("49×46=2254", $true, $false, $false, $false, $false, $true, 1, $false, "36×13=468", 2) | Out-Null
$d.Content.Find.Execute("95×95=9025", $true, $false, $false, $false, $false, $true, 1, $false, "43×71=3053", 2) | Out-Null
$d.Content.Find.Execute("17×47=799", $true, $false, $false, $false, $false, $true, 1, $false, "77×48=3696", 2) | Out-Null
$d.Content.Find.Execute("49×47=2303", $true, $false, $false, $false, $false, $true, 1, $false, "89×99=8811", 2) | Out-Null
$d.Content.Find.Execute("79×67=5293", $true, $false, $false, $false, $false, $true, 1, $false, "97×17=1649", 2) | Out-Null
